# Atualização automática via Streamlit (13/11/2025 19:00)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRINCIPAL")

$ws.Range("B3").Value = "00p0012"
$ws.Range("F3").Value = "k"
$ws.Range("G3").Value = "k"
$ws.Range("H3").Value = "k - (k 01/10/25_24h) - AM"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "01/10/25"
$ws.Range("J3").Value = "24h"
